$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 619
    $ws.Range("F3").Value = 472
    $ws.Range("F6").Value = 44
    $ws.Range("F7").Value = 37
    $ws.Range("F8").Value = 1118
    $ws.Range("F9").Value = 3887
}
